# Generate Report for Handback
#
# - Marks the "Status" cells (Overview E/F, zh-cn/de-de "Status" col C) as
#   handed back & in sync with en-US.
# - Widens the now-longer Status columns, and the "Latest Handback File"
#   column (which will hold a long .xlf filename).
# - Fills in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns for both locales now that the
#   handback xliffs have been generated, and hyperlinks the new target
#   file cell back to the source markdown file (same as column A).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5be971ca8e997e93ba2ddcb913a1f6a1491f8fd3/e2e/a.md"

# Column width, in "character" units, that round-trips to the new wider
# stored width used for the Status / Latest Handback File columns.
$wideWidth = 29.16666666667
$maxWidth = 39.16666666667

# ---------------------------------------------------------------------
# Overview sheet: update the status text shown for both locales and
# widen the two locale-status columns (E = zh-cn, F = de-de).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Range("E:E").ColumnWidth = $wideWidth
$overview.Range("F:F").ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# zh-cn sheet: status text, widen Status + Latest Handback File columns,
# and record the handback details now that the zh-cn xliff has round
# tripped back in sync with the source.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("C:C").ColumnWidth = $wideWidth
$zhcn.Range("J:J").ColumnWidth = $maxWidth

$zhcn.Range("I2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aUrl, "", "", "a.md")
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 18:42:10"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aUrl, "", "", "a.md")
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-05 18:42:10"

# ---------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, with its own xliff name
# and handback timestamp.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText
$dede.Range("C:C").ColumnWidth = $wideWidth
$dede.Range("J:J").ColumnWidth = $maxWidth

$dede.Range("I2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("I2"), $aUrl, "", "", "a.md")
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 18:42:18"

$dede.Range("I3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("I3"), $aUrl, "", "", "a.md")
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-09-05 18:42:18"
